# Disable smart-quote / autoformat substitutions so straight apostrophes in the
# replacement text survive verbatim (matches target OOXML which uses U+0027).
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [int]$ParaIndex,
        [string]$Old,
        [string]$New
    )
    $p = $d.Paragraphs($ParaIndex)
    $pStart = $p.Range.Start
    $text = $p.Range.Text
    $idx = $text.IndexOf($Old)
    if ($idx -lt 0) {
        throw "Replace-Text: substring not found in paragraph $ParaIndex : $Old"
    }
    $s = $pStart + $idx
    $e = $s + $Old.Length
    $rng = $d.Range($s, $e)
    $rng.Text = $New
}

# --- Paragraph 1: title ---------------------------------------------------
Replace-Text 1 `
    "Genomics - Unveiling the Blueprint of Life" `
    "The Profound Impact of Mathematics: Unlocking the World's Intricate Patterns"

# --- Paragraph 2: author name ----------------------------------------------
Replace-Text 2 `
    "Dr. Eleanor Hayes" `
    "Oliver Bard"

# --- Paragraph 3: email ----------------------------------------------------
Replace-Text 3 `
    "ehayes@crimsonlabs" `
    "oliver.bard98@eduworld"

# --- Paragraph 5: body (apply right-to-left so earlier offsets stay valid) -
Replace-Text 5 `
    " Genomics has also empowered biotechnology, providing the tools to harness microorganisms for the production of medicines, biofuels, and materials, unlocking the potential for a more sustainable and environmentally friendly future" `
    " Through mathematical exploration, we hone our ability to analyze, abstract, and synthesize information. Mathematics transcends cultural and linguistic boundaries, fostering a universal language for exploration and communication"

Replace-Text 5 `
    " The decipherment of genomes aids in unraveling complex ecological interactions, guiding conservation efforts and safeguarding biodiversity" `
    " It cultivates logical reasoning, problem-solving skills, and creative thinking"

Replace-Text 5 `
    "Moreover, genomics has revolutionized agriculture, enabling scientists to engineer crops resistant to pests, tolerant to changing climates, and enriched with essential nutrients" `
    "Mathematics offers a unique way of thinking that extends beyond mere computation"

Replace-Text 5 `
    " Studying genomes unveils mysteries of disease susceptibility, drug responses, and the origins of variations that define each unique individual" `
    " Its inherent beauty, logic, and practicality allure countless individuals to delve into its depths"

Replace-Text 5 `
    " These molecules drive cellular processes, shape organisms' features, and orchestrate the intricate dance of life" `
    " Mathematics underlies our understanding of everything from celestial mechanics to the behavior of subatomic particles"

Replace-Text 5 `
    " The genome, an elegant symphony of nucleotides, orchestrates the synthesis of proteins - the building blocks of living structures" `
    " These principles enable us to make sense of natural phenomena, predict outcomes, and derive insights from vast amounts of data"

Replace-Text 5 `
    "Delving into the microscopic universe of genomes grants researchers an unprecedented window into the very essence of life" `
    "The world we inhabit is governed by mathematical principles"

Replace-Text 5 `
    " By deciphering these genetic texts, genomics reveals the tale of evolution, the inheritance of traits, and the predisposition to diseases. The pursuit of genomics has revolutionized medicine, agriculture, and biotechnology, promising transformative insights into human health, the functioning of ecosystems, and the manipulation of life itself" `
    " From intricate calculations to elegant formulas, mathematics plays a pivotal role in comprehending the universe around us"

Replace-Text 5 `
    " It explores the complete set of genetic instructions, the genome, which holds the blueprints for the development, function, and replication of all living organisms" `
    " As we embark on this intellectual journey, we will unravel the profound impact of mathematics in shaping our world and transforming our perspectives"

Replace-Text 5 `
    "Genomics, a field at the frontiers of modern science, unveils the intricate tapestry of genetic information that governs life's myriad intricacies" `
    "Mathematics, a language of numbers, patterns, and structures, has served as a fundamental pillar of human understanding for millennia"

# --- Paragraph 7: summary ---------------------------------------------------
Replace-Text 7 `
    " It has empowered biotechnology industries, unlocking the potential of microorganisms for sustainable solutions. As our understanding of genomes deepens, we delve closer to the heart of life's mysteries, forging a path towards a healthier, more sustainable future" `
    " The study of mathematics goes beyond mere symbol manipulation; it is a journey of discovery, where we unravel the elegance and beauty of hidden patterns and relationships, shaping our understanding of the cosmos"

Replace-Text 7 `
    " Genomics has revolutionized agriculture, breeding crops that can withstand the challenges of a changing climate" `
    " Not only does it provide essential tools for problem-solving and decision-making, but it also cultivates analytical thinking, creativity, and a deeper appreciation for the world around us"

Replace-Text 7 `
    " By unraveling genetic instructions, it reveals the story of evolution, unravels disease mechanisms, and directs innovative medical interventions" `
    " Mathematics serves as a gateway to unlocking the intricate patterns and hidden relationships that govern the universe"

Replace-Text 7 `
    "Genomics, the study of genomes, holds the key to deciphering the intricate symphony of life" `
    "In this essay, we explored the transformative power of mathematics in our lives"

# --- Add a trailing empty paragraph after the summary paragraph ------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output "edits applied"
